$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 108, pushing the existing rows 108:120 down to 109:121
$ws.Rows(108).Insert()

# Populate the newly inserted row 108 with the latest weekly observation
$ws.Cells.Item(108, 1).Value = 5
$ws.Cells.Item(108, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(108, 3).Value = "Maule"
$ws.Cells.Item(108, 4).Value = 45132
$ws.Cells.Item(108, 5).Value = 7
$ws.Cells.Item(108, 6).Value = 100112013
$ws.Cells.Item(108, 7).Value = "Alcachofa"
$ws.Cells.Item(108, 8).Value = "Madrigal"
$ws.Cells.Item(108, 9).Value = "Primera"
$ws.Cells.Item(108, 10).Value = 300
$ws.Cells.Item(108, 11).Value = 15000
$ws.Cells.Item(108, 12).Value = 15000
$ws.Cells.Item(108, 13).Value = 15000
$ws.Cells.Item(108, 14).Value = "`$/caja 40 unidades"
$ws.Cells.Item(108, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(108, 16).Value = 375
$ws.Cells.Item(108, 17).Value = 40
$ws.Cells.Item(108, 18).Value = "Hortaliza"
